# Fix dply issue in Tourism goal trend calculations and add the new
# "T" (Tourism) goal data block to the worksheet, updating documentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Tourism" goal row (mirrors the existing "Eco"/"Liv" blocks above it):
# column A holds the goal code, C the rank, D the score, E the trend.
$ws.Range("A10").Value2 = "T"
$ws.Range("C10").Value2 = 1
$ws.Range("D10").Value2 = 62.5
$ws.Range("E10").Value2 = -0.29

$ws.Range("C11").Value2 = 2
$ws.Range("D11").Value2 = 65.2
$ws.Range("E11").Value2 = -0.28

$ws.Range("C12").Value2 = 3
$ws.Range("D12").Value2 = 59.1
$ws.Range("E12").Value2 = -0.38

$ws.Range("C13").Value2 = 4
$ws.Range("D13").Value2 = 62.4
$ws.Range("E13").Value2 = -0.38

# Move the active selection the way it was left after entering the table.
$ws.Range("A18").Select()
